$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "Мужчины"
$ws.Range("B8").Value = "Женщины"
$ws.Range("A7").Value = "Эркектер"
$ws.Range("A8").Value = "Аялдар"
$ws.Range("B6").Value = "По полу"
$ws.Range("A6").Value = "Жынысы боюнча"
$ws.Range("A10").Value = "Шаар"
$ws.Range("A11").Value = "Айыл"
$ws.Range("C6").Value = "By sex"
$ws.Range("C7").Value = "Men"
$ws.Range("C8").Value = "Woman"
$ws.Range("C10").Value = "Urban"
$ws.Range("C11").Value = "Rural"
$ws.Range("C28").Value = "Does not attend"
$ws.Range("C30").Value = "Preschool or not /primary"
$ws.Range("C31").Value = "Basic general"
$ws.Range("C32").Value = "Average total"
$ws.Range("C33").Value = "Vocational primary /secondary"
$ws.Range("C34").Value = "Higher"
$ws.Range("A35").Value = "Баланын функционалдык кыйнчылыктары"
$ws.Range("C38").Value = "Wealth quintile"
$ws.Range("C29").Value = "Educationof mother"
